$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 / Row 14 swap (ranking order changed: Chainlink now above WrappedEther) ---
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"

# --- Row 36 / Row 37 swap (ranking order changed: ImmutableX now above TrustWalletToken) ---
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"

# --- Price (column D) and Volume(1h) (column E) updates ---
$ws.Range("D2").Value = "35.092.97"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "1.850.33"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.61"
$ws.Range("E5").Value = "  +3.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.54"
$ws.Range("E8").Value = "  +6.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.327"
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0692"
$ws.Range("E10").Value = "  +1.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0989"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "2.116.64"
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.40"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("D14").Value = "1.847.14"
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.675"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.75"
$ws.Range("E16").Value = "  +3.25%  "
$ws.Range("D17").Value = "35.020.07"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.99"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").Value = "0.0₃0792"
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.53"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.16"
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.75"
$ws.Range("E22").Value = "  +1.70%  "
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "169.56"
$ws.Range("E25").Value = "  -2.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.01"
$ws.Range("E26").Value = "  +3.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.84"
$ws.Range("E27").Value = "  +21.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.60"
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.124"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0553"
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.99"
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.02"
$ws.Range("E33").Value = "  +1.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.71"
$ws.Range("E34").Value = "  +25.31%  "
$ws.Range("E35").Value = "  +9.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.797"
$ws.Range("E36").Value = "  +15.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.30"
$ws.Range("E37").Value = "  +3.66%  "
$ws.Range("E38").Value = "  +10.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0201"
$ws.Range("E39").Value = "  +4.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "90.19"
$ws.Range("E40").Value = "  -1.73%  "
$ws.Range("D41").Value = "1.342.75"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.39"
$ws.Range("E42").Value = "  +57.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.89"
$ws.Range("E43").Value = "  +3.30%  "
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0551"
$ws.Range("E47").Value = "  +5.93%  "
$ws.Range("E48").Value = "  +4.41%  "
$ws.Range("D49").Value = "2.039.72"
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0675"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("E51").Value = "  +0.28%  "
